# Edit script: rename sheet1, update numeric result values on sheets 1-4,
# and add a new "REZK" worksheet (copied structure from "BrAD") with its own data.

$wb = $excel.ActiveWorkbook

# 1. Rename "20_ISIC_DiDI" -> "20-ISIC-DiDI"
$ws1 = $wb.Worksheets.Item("20_ISIC_DiDI")
$ws1.Name = "20-ISIC-DiDI"

$ws2 = $wb.Worksheets.Item("ISIC_DiDI")
$ws3 = $wb.Worksheets.Item("DiDI")
$ws4 = $wb.Worksheets.Item("BrAD")

# 2. Update the metric values on the four existing sheets.

# Sheet 1: 20-ISIC-DiDI
$ws1.Range("B2").Value = 0.82890625
$ws1.Range("C2").Value = 0.8902795167807954
$ws1.Range("D2").Value = 0.7857095047591347
$ws1.Range("E2").Value = 0.9067708333333333
$ws1.Range("F2").Value = 0.7857095047591347
$ws1.Range("G2").Value = 0.8414532847882495
$ws1.Range("B3").Value = 0.9567708333333333
$ws1.Range("C3").Value = 0.9922172464298219
$ws1.Range("D3").Value = 0.9266042355580453
$ws1.Range("E3").Value = 0.9927083333333333
$ws1.Range("F3").Value = 0.9266042355580453
$ws1.Range("G3").Value = 0.9583827962286656
$ws1.Range("B4").Value = 0.8580729166666666
$ws1.Range("C4").Value = 0.983363499052815
$ws1.Range("D4").Value = 0.7851087919200189
$ws1.Range("E4").Value = 0.9875
$ws1.Range("F4").Value = 0.7851087919200189
$ws1.Range("G4").Value = 0.8745400347457845
$ws1.Range("B5").Value = 0.8450520833333334
$ws1.Range("C5").Value = 0.92221541952537
$ws1.Range("D5").Value = 0.7928836980892566
$ws1.Range("E5").Value = 0.9354166666666667
$ws1.Range("F5").Value = 0.7928836980892566
$ws1.Range("G5").Value = 0.8579287640738509
$ws1.Range("B6").Value = 0.9502604166666667
$ws1.Range("C6").Value = 0.9932736999445972
$ws1.Range("D6").Value = 0.9147769671847101
$ws1.Range("E6").Value = 0.99375
$ws1.Range("F6").Value = 0.9147769671847101
$ws1.Range("G6").Value = 0.9524653153464017
$ws1.Range("B7").Value = 0.8348958333333333
$ws1.Range("C7").Value = 0.9158589165885606
$ws1.Range("D7").Value = 0.7809402427834973
$ws1.Range("E7").Value = 0.9317708333333333
$ws1.Range("F7").Value = 0.7809402427834973
$ws1.Range("G7").Value = 0.8495253283407724

# Sheet 2: ISIC_DiDI
$ws2.Range("B2").Value = 0.85078125
$ws2.Range("C2").Value = 0.8899047454384987
$ws2.Range("D2").Value = 0.8214462830222615
$ws2.Range("E2").Value = 0.8994791666666667
$ws2.Range("F2").Value = 0.8214462830222615
$ws2.Range("G2").Value = 0.8579073260141875
$ws2.Range("B3").Value = 0.96953125
$ws2.Range("C3").Value = 0.9866226778404528
$ws2.Range("D3").Value = 0.9540059264886415
$ws2.Range("E3").Value = 0.9869791666666666
$ws2.Range("F3").Value = 0.9540059264886415
$ws2.Range("G3").Value = 0.970108220629583
$ws2.Range("B4").Value = 0.8729166666666667
$ws2.Range("C4").Value = 0.9678422940670781
$ws2.Range("D4").Value = 0.8116010939347308
$ws2.Range("E4").Value = 0.9739583333333334
$ws2.Range("F4").Value = 0.8116010939347308
$ws2.Range("G4").Value = 0.884953766421864
$ws2.Range("B5").Value = 0.865625
$ws2.Range("C5").Value = 0.933628230533488
$ws2.Range("D5").Value = 0.8180479111908884
$ws2.Range("E5").Value = 0.94375
$ws2.Range("F5").Value = 0.8180479111908884
$ws2.Range("G5").Value = 0.8758230140092156
$ws2.Range("B6").Value = 0.9671875
$ws2.Range("C6").Value = 0.9854259227781834
$ws2.Range("D6").Value = 0.9506203236695991
$ws2.Range("E6").Value = 0.9859375
$ws2.Range("F6").Value = 0.9506203236695991
$ws2.Range("G6").Value = 0.9678662161876035
$ws2.Range("B7").Value = 0.8526041666666667
$ws2.Range("C7").Value = 0.9293606508088482
$ws2.Range("D7").Value = 0.8005391311694824
$ws2.Range("E7").Value = 0.9411458333333333
$ws2.Range("F7").Value = 0.8005391311694824
$ws2.Range("G7").Value = 0.8647474916355496

# Sheet 3: DiDI
$ws3.Range("B2").Value = 0.5010416666666667
$ws3.Range("C2").Value = 0.5020401497205079
$ws3.Range("D2").Value = 0.5007026800932429
$ws3.Range("E2").Value = 0.7208333333333333
$ws3.Range("F2").Value = 0.5007026800932429
$ws3.Range("G2").Value = 0.5907457510130715
$ws3.Range("B3").Value = 0.959375
$ws3.Range("C3").Value = 0.9912511499384221
$ws3.Range("D3").Value = 0.9318695247326225
$ws3.Range("E3").Value = 0.9916666666666667
$ws3.Range("F3").Value = 0.9318695247326225
$ws3.Range("G3").Value = 0.9607112127740762
$ws3.Range("B4").Value = 0.5411458333333333
$ws3.Range("C4").Value = 0.6981170622463193
$ws3.Range("D4").Value = 0.523050500505784
$ws3.Range("E4").Value = 0.9333333333333333
$ws3.Range("F4").Value = 0.523050500505784
$ws3.Range("G4").Value = 0.6703386882164902
$ws3.Range("B5").Value = 0.534375
$ws3.Range("C5").Value = 0.5765415166730162
$ws3.Range("D5").Value = 0.5223352010449822
$ws3.Range("E5").Value = 0.7973958333333333
$ws3.Range("F5").Value = 0.5223352010449822
$ws3.Range("G5").Value = 0.6309781514637163
$ws3.Range("B6").Value = 0.9682291666666667
$ws3.Range("C6").Value = 0.9946113259294572
$ws3.Range("D6").Value = 0.9451248822817467
$ws3.Range("E6").Value = 0.9947916666666666
$ws3.Range("F6").Value = 0.9451248822817467
$ws3.Range("G6").Value = 0.969166687676167
$ws3.Range("B7").Value = 0.4997395833333333
$ws3.Range("C7").Value = 0.4999098482288729
$ws3.Range("D7").Value = 0.499815256527091
$ws3.Range("E7").Value = 0.8208333333333333
$ws3.Range("F7").Value = 0.499815256527091
$ws3.Range("G7").Value = 0.6211640040857177

# Sheet 4: BrAD
$ws4.Range("B2").Value = 0.5489583333333333
$ws4.Range("C2").Value = 0.5891536835158014
$ws4.Range("D2").Value = 0.5338912074759978
$ws4.Range("E2").Value = 0.7890625
$ws4.Range("F2").Value = 0.5338912074759978
$ws4.Range("G2").Value = 0.6363516621066385
$ws4.Range("B3").Value = 0.5598958333333334
$ws4.Range("C3").Value = 0.619242440019829
$ws4.Range("D3").Value = 0.5403948433656248
$ws4.Range("E3").Value = 0.809375
$ws4.Range("F3").Value = 0.5403948433656248
$ws4.Range("G3").Value = 0.647712764775512
$ws4.Range("B4").Value = 0.5075520833333333
$ws4.Range("C4").Value = 0.5317490961442006
$ws4.Range("D4").Value = 0.5042950462172381
$ws4.Range("E4").Value = 0.8463541666666666
$ws4.Range("F4").Value = 0.5042950462172381
$ws4.Range("G4").Value = 0.6317916588874618
$ws4.Range("B5").Value = 0.53359375
$ws4.Range("C5").Value = 0.5714667078242347
$ws4.Range("D5").Value = 0.5221417460128234
$ws4.Range("E5").Value = 0.7979166666666667
$ws4.Range("F5").Value = 0.5221417460128234
$ws4.Range("G5").Value = 0.6309355869408093
$ws4.Range("B6").Value = 0.5559895833333334
$ws4.Range("C6").Value = 0.6099966214931587
$ws4.Range("D6").Value = 0.5378592576596339
$ws4.Range("E6").Value = 0.7947916666666667
$ws4.Range("F6").Value = 0.5378592576596339
$ws4.Range("G6").Value = 0.6412538169887511
$ws4.Range("B7").Value = 0.540625
$ws4.Range("C7").Value = 0.6190579031213258
$ws4.Range("D7").Value = 0.5246495308745153
$ws4.Range("E7").Value = 0.8755208333333333
$ws4.Range("F7").Value = 0.5246495308745153
$ws4.Range("G7").Value = 0.655960186406647

# 3. Add the new "REZK" sheet at the end, copying the layout/styles of "BrAD",
#    then overwrite its values with the new dataset's results.
$ws4.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "REZK"

# Overwrite REZK data
$ws5.Range("A2").Value = 0
$ws5.Range("B2").Value = 0.5526041666666667
$ws5.Range("C2").Value = 0.6549639690325181
$ws5.Range("D2").Value = 0.5317837525038128
$ws5.Range("E2").Value = 0.8921875
$ws5.Range("F2").Value = 0.5317837525038128
$ws5.Range("G2").Value = 0.6662450133159682
$ws5.Range("A3").Value = 1
$ws5.Range("B3").Value = 0.5653645833333333
$ws5.Range("C3").Value = 0.6671130619485622
$ws5.Range("D3").Value = 0.5410074393019572
$ws5.Range("E3").Value = 0.8744791666666667
$ws5.Range("F3").Value = 0.5410074393019572
$ws5.Range("G3").Value = 0.6681982559939662
$ws5.Range("A4").Value = 2
$ws5.Range("B4").Value = 0.54296875
$ws5.Range("C4").Value = 0.7535540385540387
$ws5.Range("D4").Value = 0.5235627476027453
$ws5.Range("E4").Value = 0.9604166666666667
$ws5.Range("F4").Value = 0.5235627476027453
$ws5.Range("G4").Value = 0.6776363306495852
$ws5.Range("A5").Value = 3
$ws5.Range("B5").Value = 0.55859375
$ws5.Range("C5").Value = 0.8420996133496134
$ws5.Range("D5").Value = 0.5322199886836917
$ws5.Range("E5").Value = 0.9729166666666667
$ws5.Range("F5").Value = 0.5322199886836917
$ws5.Range("G5").Value = 0.6879816891958406
$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = 0.55859375
$ws5.Range("C6").Value = 0.6611588083908793
$ws5.Range("D6").Value = 0.5360671837765901
$ws5.Range("E6").Value = 0.8786458333333333
$ws5.Range("F6").Value = 0.5360671837765901
$ws5.Range("G6").Value = 0.6656759121339679
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = 0.4971354166666667
$ws5.Range("C7").Value = 0.4852966795783796
$ws5.Range("D7").Value = 0.4984546480210013
$ws5.Range("E7").Value = 0.8536458333333333
$ws5.Range("F7").Value = 0.4984546480210013
$ws5.Range("G7").Value = 0.6293188090104486

$ws1.Select()
